{"js": "// Q2: \": Which two land-use classes constitute the most flood damage? Which percentage...\"\n// becomes\n// \": Which two land-use classes constitute the most flood damage for a 1/100 flood event? Which percentage...\"\nconst q2Results = context.document.body.search(\"most flood damage\", { matchCase: true });\nq2Results.load(\"items\");\nawait context.sync();\n\nif (q2Results.items.length > 0) {\n  q2Results.items[0].insertText(\" for a 1/100 flood event\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Q3: \" Create a figure that shows the distribution of losses. You can either save the \"\n// becomes\n// \" Create a figure that shows the distribution of losses for a 1/100 flood event. You can either save the \"\n// (with a collapsed \"_GoBack\" bookmark right before \". You can either save the\")\nconst q3Results = context.document.body.search(\"distribution of losses\", { matchCase: true });\nq3Results.load(\"items\");\nawait context.sync();\n\nif (q3Results.items.length > 0) {\n  q3Results.items[0].insertText(\" for a 1/100 flood event\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Re-locate the text we just inserted in Q3 (unique in the document) so the\n// \"_GoBack\" bookmark can be dropped as a collapsed point right after it.\nconst q3Inserted = context.document.body.search(\"distribution of losses for a 1/100 flood event\", { matchCase: true });\nq3Inserted.load(\"items\");\nawait context.sync();\n\nif (q3Inserted.items.length > 0) {\n  const afterPoint = q3Inserted.items[0].getRange(Word.RangeLocation.after);\n  afterPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Q2: \": Which two land-use classes constitute the most flood damage? Which percentage...\"\n# becomes\n# \": Which two land-use classes constitute the most flood damage for a 1/100 flood event? Which percentage...\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"most flood damage\"\n$find.MatchCase = $true\n$found = $find.Execute()\nif ($found) {\n    $rng = $find.Parent.Duplicate\n    $rng.Collapse(0)\n    $rng.InsertAfter(\" for a 1/100 flood event\")\n}\n\n# Q3: \" Create a figure that shows the distribution of losses. You can either save the \"\n# becomes\n# \" Create a figure that shows the distribution of losses for a 1/100 flood event. You can either save the \"\n# (with a collapsed \"_GoBack\" bookmark right before \". You can either save the\")\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"distribution of losses\"\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\nif ($found2) {\n    $rng2 = $find2.Parent.Duplicate\n    $rng2.Collapse(0)\n    $rng2.InsertAfter(\" for a 1/100 flood event\")\n    $rng2.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $rng2)\n}\n"}
